# Actualización del Contenido de los Juegos
#
# 1) Swap the two embedded-font declarations (Calibri <-> Raleway) in the
#    presentation's embeddedFontLst (the embedded font *data* r:id links
#    stay attached to the same <p:embeddedFont> blocks; only the
#    <p:font .../> identity each block declares is swapped).
# 2) On slide 2, shape 1 ("- 7 tarjetas donde se explican ...") split the
#    leading "- 7 " into its own run with updated text "- 8 ", leaving the
#    remainder of the sentence in a second run.

$p = $ppt.ActivePresentation

# --- 1) Swap embedded font declarations -----------------------------------
$fonts = $p.EmbeddedFontLst
if ($fonts -ne $null) {
    $f1 = $fonts.Item(1)
    $f3 = $fonts.Item(3)
    $tmpTypeface = $f1.Font.Name
    $name1 = $f1.Font.Name
    $name3 = $f3.Font.Name
    $f1.Font.Name = $name3
    $f3.Font.Name = $name1
}

# --- 2) Update the "- 7 tarjetas" text -------------------------------------
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

$target = $tr.Find("- 7 ")
$target.Text = "- 8 "
